# "change wechat 2 QQ group"
#
# The "Resources" slide (slide 12) contains a SmartArt diagram
# (graphicFrame "Content Placeholder 3") with a node whose text reads
# "WeChat Group: programming windows 2020". The author replaced the
# WeChat mention with a QQ group number and the two students' course
# citations ("20201021076<Windows原理与应用>" / "20201021952<...>").
#
# We locate that SmartArt node via Shape.SmartArt.AllNodes and rewrite
# its text in place; PowerPoint keeps the diagram data (dgm:t) and the
# cached diagram drawing (dsp:txBody) in sync automatically.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$sh = $s.Shapes.Item("Content Placeholder 3")

$sa = $sh.SmartArt
$nodes = $sa.AllNodes

for ($i = 1; $i -le $nodes.Count; $i++) {
    $node = $nodes.Item($i)
    $tr = $node.TextFrame2.TextRange
    if ($tr.Text -like "WeChat Group:*") {
        $tr.Text = "QQ group: `r20201021076《Windows原理与应用》`r20201021952《Windows原理与应用》"
        break
    }
}
